$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 80
    3  = 3094
    5  = 2603
    8  = 6
    9  = 1339
    13 = 1173
    14 = 342
    15 = 323
    16 = 32
    20 = 87
    21 = 2431
    23 = 277
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
